$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so that values
# like "327.49" or "1.002" are not auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.218.76"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "1.910.17"
$ws.Range("E3").Value = "  +1.82%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").Value = "327.49"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "0.4617"
$ws.Range("E7").Value = "  +0.57%  "

$ws.Range("D8").Value = "0.3944"
$ws.Range("E8").Value = "  +1.88%  "

$ws.Range("D9").Value = "46.77"
$ws.Range("E9").Value = "  +1.41%  "

$ws.Range("D10").Value = "0.07944"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  +1.27%  "

$ws.Range("D12").Value = "22.33"
$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("D13").Value = "1.926.61"
$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("D14").Value = "7.108"
$ws.Range("E14").Value = "  +1.58%  "

$ws.Range("D15").Value = "5.768"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "0.06950"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "88.47"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").Value = "0.00001008"

$ws.Range("D20").Value = "17.15"

$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "29.234.79"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("D23").Value = "5.366"
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").Value = "2.147.35"
$ws.Range("E25").Value = "  +1.15%  "

$ws.Range("D26").Value = "2.058"
$ws.Range("E26").Value = "  -2.17%  "

$ws.Range("D27").Value = "156.87"
$ws.Range("E27").Value = "  +2.41%  "

$ws.Range("D28").Value = "19.50"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").Value = "6.125"
$ws.Range("E29").Value = "  +5.52%  "

$ws.Range("D30").Value = "1.999"
$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("D31").Value = "118.82"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").Value = "0.09378"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").Value = "0.9287"
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("D34").Value = "5.353"
$ws.Range("E34").Value = "  +0.82%  "

$ws.Range("D35").Value = "1.353"
$ws.Range("E35").Value = "  +0.88%  "

$ws.Range("D36").Value = "3.273"
$ws.Range("E36").Value = "  -1.54%  "

$ws.Range("D37").Value = "1.208"
$ws.Range("E37").Value = "  +4.92%  "

$ws.Range("D38").Value = "0.05836"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").Value = "0.02108"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").Value = "7.934"
$ws.Range("E40").Value = "  +3.05%  "

$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("D42").Value = "0.5759"
$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("D43").Value = "0.1803"

$ws.Range("D44").Value = "9.979"
$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("D45").Value = "2.259"
$ws.Range("E45").Value = "  +4.41%  "

$ws.Range("D46").Value = "11.93"
$ws.Range("E46").Value = "  +0.74%  "

$ws.Range("D47").Value = "0.5414"
$ws.Range("E47").Value = "  +2.15%  "

$ws.Range("D48").Value = "0.07075"
$ws.Range("E48").Value = "  -2.03%  "

$ws.Range("D49").Value = "1.878"
$ws.Range("E49").Value = "  +2.75%  "

$ws.Range("D50").Value = "2.555"
$ws.Range("E50").Value = "  +6.08%  "

$ws.Range("D51").Value = "113.05"
$ws.Range("E51").Value = "  -0.69%  "

# Restore the original (default/no explicit) cell style for column D
# now that the text values are safely stored, so no stray style index
# is left referenced on these cells.
$ws.Range("D2:D51").Style = "Normal"
